$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the three runs that spell out the "master" branch rename bullet
#    ("To Rename a \u201c" + "m" + "aster\u201d  branch to \"main\", use.")
#    into a single run with identical text, by re-writing just the leading
#    "To Rename a" text in place - Word collapses the touched run together
#    with its immediately adjacent same-formatted runs.
# ---------------------------------------------------------------------------
$renameRange = $d.Content
$renameRange.Find.Execute("To Rename a", $true, $false, $false, $false, $false, $true, 1, $false, "To Rename a", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Remove the "Pushing the stagged files to GitHub Repo." / "Try It:" /
#    "git push origin main" block - this content was pushed into a new
#    branch per the commit message and is no longer needed on this page.
#    The block spans from the (empty, Consolas-formatted) paragraph right
#    before "Pushing the stagged files..." through two empty
#    Consolas-formatted paragraphs that trail "git push origin main".
# ---------------------------------------------------------------------------
$pushFind = $d.Content
$pushFind.Find.Execute("Pushing the stagged files to GitHub Repo.") | Out-Null
$startIdx = $pushFind.Paragraphs(1).Index

$mainFind = $d.Content
$mainFind.Find.Execute("push origin main") | Out-Null
$endIdx = $mainFind.Paragraphs(1).Index

$deleteStart = $d.Paragraphs($startIdx - 1).Range.Start
$deleteEnd = $d.Paragraphs($endIdx + 2).Range.End
$d.Range($deleteStart, $deleteEnd).Delete()

# ---------------------------------------------------------------------------
# 3) Drop the stray empty paragraph sitting right after the bookmarked
#    (Toc) paragraph near the top of the "Configure Git" section.
# ---------------------------------------------------------------------------
$emailFind = $d.Content
$emailFind.Find.Execute("set an email address that will be associated with each history marker") | Out-Null
$emailIdx = $emailFind.Paragraphs(1).Index
# emailIdx     -> "...history marker" paragraph
# emailIdx + 1 -> bookmark paragraph (empty)
# emailIdx + 2 -> stray empty paragraph to remove
$d.Paragraphs($emailIdx + 2).Range.Delete()

Write-Host "Paragraphs now:" $d.Paragraphs.Count
